$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old trailing rows (31-49) entirely, preserving formatting of remaining rows
$ws.Range("A31:B49").EntireRow.Delete()

# Header
$ws.Range("B1").Value = "urls"

# Data rows 2-30: refreshed product URL list
$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(2, 2).Value = "https://www.uyyaala.com/products/abbott-similac-gold-stage-1-new-advanced-milk-formula-with-hmo-400g-0-6m"
$ws.Cells.Item(3, 1).Value = 1
$ws.Cells.Item(3, 2).Value = "https://www.uyyaala.com/products/abbott-similac-gold-stage-3-new-advanced-milk-formula-with-hmo-400g-1-3y"
$ws.Cells.Item(4, 1).Value = 2
$ws.Cells.Item(4, 2).Value = "https://www.uyyaala.com/products/abbott-similac-gold-stage-2-new-advanced-milk-formula-with-hmo-400g-6-to-12m"
$ws.Cells.Item(5, 1).Value = 3
$ws.Cells.Item(5, 2).Value = "https://www.uyyaala.com/products/abbott-similac-infant-formula-stage-1-1-to-6-months-400g"
$ws.Cells.Item(6, 1).Value = 4
$ws.Cells.Item(6, 2).Value = "https://www.uyyaala.com/products/abbott-pediasure-complete-and-balance-nutrition-vanilla-1-10-years-400-gmimported"
$ws.Cells.Item(7, 1).Value = 5
$ws.Cells.Item(7, 2).Value = "https://www.uyyaala.com/products/abbott-similac-advance-2-after-6-months"
$ws.Cells.Item(8, 1).Value = 6
$ws.Cells.Item(8, 2).Value = "https://www.uyyaala.com/products/abbott-similac-advance-1-up-to-6-months"
$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(9, 2).Value = "https://www.uyyaala.com/products/abbott-similac-iq-infant-formula-stage-1-0-to-6-months-400g"
$ws.Cells.Item(10, 1).Value = 8
$ws.Cells.Item(10, 2).Value = "https://www.uyyaala.com/products/abbott-similac-infant-formula-stage-1-400g-tin-0-6m"
$ws.Cells.Item(11, 1).Value = 9
$ws.Cells.Item(11, 2).Value = "https://www.uyyaala.com/products/abbott-pediasure-complete-triple-sure-strawberry-1-10-years-400gimported"
$ws.Cells.Item(12, 1).Value = 10
$ws.Cells.Item(12, 2).Value = "https://www.uyyaala.com/products/abbott-similac-isomil-lactose-free-infant-milk-substitute-0-to-24-months-400g"
$ws.Cells.Item(13, 1).Value = 11
$ws.Cells.Item(13, 2).Value = "https://www.uyyaala.com/products/abbott-pediasure-complete-and-balance-nutrition-chocolate-1-10-years-400-gmimported"
$ws.Cells.Item(14, 1).Value = 12
$ws.Cells.Item(14, 2).Value = "https://www.uyyaala.com/products/abbott-similac-follow-up-formula-stage-2-after-6-months"
$ws.Cells.Item(15, 1).Value = 13
$ws.Cells.Item(15, 2).Value = "https://www.uyyaala.com/products/abbott-similac-iq-follow-up-formula-stage-3-12-to-24-months-400g"
$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(16, 2).Value = "https://www.uyyaala.com/products/abbott-pediasure-complete-triple-sure-honey-1-10-years-400g-imported"
$ws.Cells.Item(17, 1).Value = 15
$ws.Cells.Item(17, 2).Value = "https://www.uyyaala.com/products/similac-360-total-care-infant-milk-formula-1-13kg"
$ws.Cells.Item(18, 1).Value = 16
$ws.Cells.Item(18, 2).Value = "https://www.uyyaala.com/products/similac-sensitive-360-total-care-infant-milk-formula-1-13kg"
$ws.Cells.Item(19, 1).Value = 17
$ws.Cells.Item(19, 2).Value = "https://www.uyyaala.com/products/abbott-similac-total-comfort-gold-infant-milk-formula-stage-1-360gms-0-6months-imported"
$ws.Cells.Item(20, 1).Value = 18
$ws.Cells.Item(20, 2).Value = "https://www.uyyaala.com/products/abbott-similac-follow-up-formula-stage-4-18-to-24-months-400g"
$ws.Cells.Item(21, 1).Value = 19
$ws.Cells.Item(21, 2).Value = "https://www.uyyaala.com/products/abbott-pediasure-complete-balanced-nutrition-to-help-kids-grow-box-nutrition-drink-refill-pack-chocolate-flavour-2-years-1kg"
$ws.Cells.Item(22, 1).Value = 20
$ws.Cells.Item(22, 2).Value = "https://www.uyyaala.com/products/abbott-pediasure-complete-balanced-nutrition-to-help-kids-grow-box-nutrition-drink-chocolate-flavour-2-years-200-g-jar"
$ws.Cells.Item(23, 1).Value = 21
$ws.Cells.Item(23, 2).Value = "https://www.uyyaala.com/products/abbott-pediasure-vanilla-delight-flavour-2years"
$ws.Cells.Item(24, 1).Value = 22
$ws.Cells.Item(24, 2).Value = "https://www.uyyaala.com/products/abbott-pediasure-complete-balanced-nutrition-to-help-kids-grow-box-nutrition-drink-refill-pack-vanilla-delight-flavour-2-years-200-g"
$ws.Cells.Item(25, 1).Value = 23
$ws.Cells.Item(25, 2).Value = "https://www.uyyaala.com/products/abbott-pediasure-complete-balanced-nutrition-to-help-kids-grow-box-nutrition-drink-refill-pack-vanilla-delight-flavour-2-years-400-g"
$ws.Cells.Item(26, 1).Value = 24
$ws.Cells.Item(26, 2).Value = "https://www.uyyaala.com/products/abbott-pediasure-complete-balanced-nutrition-to-help-kids-grow-box-nutrition-drink-refill-pack-chocolate-flavour-2-years-200-g"
$ws.Cells.Item(27, 1).Value = 25
$ws.Cells.Item(27, 2).Value = "https://www.uyyaala.com/products/abbott-similac-total-comfort-infant-milk-formula-0-to-6months-350gms"
$ws.Cells.Item(28, 1).Value = 26
$ws.Cells.Item(28, 2).Value = "https://www.uyyaala.com/products/abbott-similac-total-comfort-infant-formula-360g-stage-2-from-6-12-months-imported"
$ws.Cells.Item(29, 1).Value = 27
$ws.Cells.Item(29, 2).Value = "https://www.uyyaala.com/products/abbott-similac-total-comfort-infant-formula-360g-stage-3-from-1-3-years-imported"
$ws.Cells.Item(30, 1).Value = 28
$ws.Cells.Item(30, 2).Value = "https://www.uyyaala.com/products/similac-advance-optigro-complete-nutrition-infant-formula-for-0-12-months-352gms"
